$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that currently sits right after
#    the "На странице видны два товара для собак" run.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Remove the whole paragraph that holds the screenshot image (the
#    <w:drawing> run with the lastRenderedPageBreak) which currently follows
#    the "Скриншот, Видео" paragraph. Remember the paragraph that follows it
#    (the empty heading paragraph) so we can re-bookmark it afterwards.
$imgIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $imgIndex = $i
    }
}

if ($imgIndex -gt 0) {
    $imgPara = $d.Paragraphs.Item($imgIndex)
    $imgPara.Range.Delete()

    # 3. Re-add the "_GoBack" bookmark, now collapsed at the end of the
    #    paragraph that took the deleted paragraph's place (the empty
    #    heading paragraph that used to come right after the image).
    $headingPara = $d.Paragraphs.Item($imgIndex)
    $d.Bookmarks.Add("_GoBack", $headingPara.Range)
}
